$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.131.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.056.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.671"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.14%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0796"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +7.71%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.836"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.055.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +31.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.151.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "76.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.127"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  +7.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0632"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0901"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").Value = "  +4.88%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  +26.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0225"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.297.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.239.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
